# ---------------------------------------------------------------------------
# Commit: "Fixed update to excel issue"
#
# 1) Rename the "Requested quantity" header on "Weekly Quantity" to
#    "Weekly_PO_Qty" and on "Monthly Trend" to "Monthly_PO_Qty".
# 2) Add a new "PO Forecast" worksheet (after the existing sheets) with
#    forecast data: ds / PO_Forecast / yhat_lower / yhat_upper.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1) Rename the "Requested quantity" headers -----------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the "PO Forecast" worksheet --------------------------------------
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "PO Forecast"

$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

# Reuse the header style already used by "Weekly Quantity"!A1:B1 (bold,
# bordered, centered) for the new sheet's header row.
$wsWeekly.Range("A1:B1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

$ws.Range("A2").Value = 45270.99999999999
$ws.Range("B2").Value = 14
$ws.Range("C2").Value = -102.6145297035914
$ws.Range("D2").Value = 118.6829139844437
$ws.Range("A3").Value = 45277.99999999999
$ws.Range("B3").Value = 16
$ws.Range("C3").Value = -92.65511391561714
$ws.Range("D3").Value = 126.7704507656066
$ws.Range("A4").Value = 45298.99999999999
$ws.Range("B4").Value = 21
$ws.Range("C4").Value = -90.10245372620062
$ws.Range("D4").Value = 129.302405779324
$ws.Range("A5").Value = 45305.99999999999
$ws.Range("B5").Value = 23
$ws.Range("C5").Value = -93.52801680126804
$ws.Range("D5").Value = 131.4941073212467
$ws.Range("A6").Value = 45312.99999999999
$ws.Range("B6").Value = 25
$ws.Range("C6").Value = -90.64557246527514
$ws.Range("D6").Value = 135.5858456187621
$ws.Range("A7").Value = 45319.99999999999
$ws.Range("B7").Value = 27
$ws.Range("C7").Value = -82.70881764019514
$ws.Range("D7").Value = 136.3432411091244
$ws.Range("A8").Value = 45326.99999999999
$ws.Range("B8").Value = 29
$ws.Range("C8").Value = -84.20348622639499
$ws.Range("D8").Value = 139.0586142773058
$ws.Range("A9").Value = 45333.99999999999
$ws.Range("B9").Value = 31
$ws.Range("C9").Value = -84.54061856704878
$ws.Range("D9").Value = 148.4497314590374
$ws.Range("A10").Value = 45354.99999999999
$ws.Range("B10").Value = 37
$ws.Range("C10").Value = -68.87280158141658
$ws.Range("D10").Value = 144.4285968584393
$ws.Range("A11").Value = 45438.99999999999
$ws.Range("B11").Value = 60
$ws.Range("C11").Value = -49.93515701888379
$ws.Range("D11").Value = 176.7322697776829
$ws.Range("A12").Value = 45480.99999999999
$ws.Range("B12").Value = 71
$ws.Range("C12").Value = -41.94081734893326
$ws.Range("D12").Value = 181.932212483116
$ws.Range("A13").Value = 45487.99999999999
$ws.Range("B13").Value = 73
$ws.Range("C13").Value = -28.76543820386568
$ws.Range("D13").Value = 183.4774788499101
$ws.Range("A14").Value = 45494.99999999999
$ws.Range("B14").Value = 75
$ws.Range("C14").Value = -33.57188634716378
$ws.Range("D14").Value = 187.0991423133664
$ws.Range("A15").Value = 45501.99999999999
$ws.Range("B15").Value = 77
$ws.Range("C15").Value = -38.13670110706391
$ws.Range("D15").Value = 196.2190940554407
$ws.Range("A16").Value = 45508.99999999999
$ws.Range("B16").Value = 79
$ws.Range("C16").Value = -34.06424359297373
$ws.Range("D16").Value = 190.0909097921583
$ws.Range("A17").Value = 45515.99999999999
$ws.Range("B17").Value = 81
$ws.Range("C17").Value = -36.42954020796371
$ws.Range("D17").Value = 195.2981584848851
$ws.Range("A18").Value = 45536.99999999999
$ws.Range("B18").Value = 86
$ws.Range("C18").Value = -24.11316042298433
$ws.Range("D18").Value = 196.0354681330795
$ws.Range("A19").Value = 45550.99999999999
$ws.Range("B19").Value = 90
$ws.Range("C19").Value = -16.21150365403111
$ws.Range("D19").Value = 199.9213304685917
$ws.Range("A20").Value = 45557.99999999999
$ws.Range("B20").Value = 92
$ws.Range("C20").Value = -15.95444174471232
$ws.Range("D20").Value = 207.0144095274597
$ws.Range("A21").Value = 45564.99999999999
$ws.Range("B21").Value = 94
$ws.Range("C21").Value = -18.35919309491492
$ws.Range("D21").Value = 199.1763487255095
$ws.Range("A22").Value = 45571.99999999999
$ws.Range("B22").Value = 96
$ws.Range("C22").Value = -20.52906580375431
$ws.Range("D22").Value = 203.6118750485004
$ws.Range("A23").Value = 45578.99999999999
$ws.Range("B23").Value = 98
$ws.Range("C23").Value = -19.95508502340115
$ws.Range("D23").Value = 208.3904539401801
$ws.Range("A24").Value = 45585.99999999999
$ws.Range("B24").Value = 100
$ws.Range("C24").Value = -9.096405516355901
$ws.Range("D24").Value = 206.8777568081636
$ws.Range("A25").Value = 45592.99999999999
$ws.Range("B25").Value = 102
$ws.Range("C25").Value = -7.855574086389268
$ws.Range("D25").Value = 216.4991517706714
$ws.Range("A26").Value = 45599.99999999999
$ws.Range("B26").Value = 103
$ws.Range("C26").Value = -11.43191115816521
$ws.Range("D26").Value = 215.5479540428815
$ws.Range("A27").Value = 45606.99999999999
$ws.Range("B27").Value = 105
$ws.Range("C27").Value = -3.392208153196185
$ws.Range("D27").Value = 222.5051110980764
$ws.Range("A28").Value = 45613.99999999999
$ws.Range("B28").Value = 107
$ws.Range("C28").Value = -6.712476934202124
$ws.Range("D28").Value = 220.6436398964288
$ws.Range("A29").Value = 45620.99999999999
$ws.Range("B29").Value = 109
$ws.Range("C29").Value = -2.720045149752456
$ws.Range("D29").Value = 210.6541863673449

# Reuse the date-column style already used by "Weekly Quantity"!A2:A3 for the
# "ds" column on the new sheet.
$wsWeekly.Range("A2:A3").Copy()
$ws.Range("A2:A29").PasteSpecial(-4122)

$excel.CutCopyMode = $false
